$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 2497
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 2497
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 7491
$ws.Cells.Item(80, 13).Value = ""
$ws.Cells.Item(80, 14).Value = -9487
$ws.Cells.Item(83, 8).Value = 2497
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 2497
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 22473
$ws.Cells.Item(83, 13).Value = ""
$ws.Cells.Item(83, 14).Value = -32457
$ws.Cells.Item(111, 8).Value = 1049.3636
$ws.Cells.Item(111, 9).Value = 834.8
$ws.Cells.Item(111, 11).Value = 2504.4
$ws.Cells.Item(111, 13).Value = 562.6000000000004
$ws.Cells.Item(125, 8).Value = 7834.9
$ws.Cells.Item(125, 9).Value = 6192.857
$ws.Cells.Item(125, 11).Value = 55735.713
$ws.Cells.Item(125, 13).Value = -53275.713
$ws.Cells.Item(137, 8).Value = 2308.75
$ws.Cells.Item(137, 9).Value = 1546
$ws.Cells.Item(137, 11).Value = 4638
$ws.Cells.Item(137, 13).Value = -2088

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4195.8105
$ws.Cells.Item(32, 9).Value = 3434.9465
$ws.Cells.Item(32, 10).Value = 25500
$ws.Cells.Item(32, 11).Value = 3434.9465
$ws.Cells.Item(32, 12).Value = 25500
$ws.Cells.Item(32, 13).Value = -3147.9465
$ws.Cells.Item(32, 14).Value = -26074
$ws.Cells.Item(61, 8).Value = 1899.5
$ws.Cells.Item(61, 9).Value = 1899
$ws.Cells.Item(61, 11).Value = 1899
$ws.Cells.Item(61, 13).Value = -1687
$ws.Cells.Item(74, 8).Value = 2380.9
$ws.Cells.Item(74, 9).Value = 2601.8
$ws.Cells.Item(74, 10).Value = 2160
$ws.Cells.Item(74, 11).Value = 2601.8
$ws.Cells.Item(74, 12).Value = 2160
$ws.Cells.Item(74, 13).Value = -1727.8
$ws.Cells.Item(74, 14).Value = -3908
$ws.Cells.Item(77, 8).Value = 2380.9
$ws.Cells.Item(77, 9).Value = 2601.8
$ws.Cells.Item(77, 10).Value = 2160
$ws.Cells.Item(77, 11).Value = 13009
$ws.Cells.Item(77, 12).Value = 10800
$ws.Cells.Item(77, 13).Value = -8641
$ws.Cells.Item(77, 14).Value = -19536
$ws.Cells.Item(97, 8).Value = 996.61536
$ws.Cells.Item(97, 9).Value = 996.0909
$ws.Cells.Item(97, 11).Value = 996.0909
$ws.Cells.Item(97, 13).Value = -500.0909
$ws.Cells.Item(122, 8).Value = 2151.2632
$ws.Cells.Item(122, 9).Value = 2220.7778
$ws.Cells.Item(122, 11).Value = 6662.3334
$ws.Cells.Item(122, 13).Value = -4212.3334
$ws.Cells.Item(132, 8).Value = 7419
$ws.Cells.Item(132, 9).Value = 7989.8335
$ws.Cells.Item(132, 11).Value = 23969.5005
$ws.Cells.Item(132, 13).Value = -21439.5005
$ws.Cells.Item(136, 8).Value = 1899.5
$ws.Cells.Item(136, 9).Value = 1899
$ws.Cells.Item(136, 11).Value = 5697
$ws.Cells.Item(136, 13).Value = -3147

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 297.5
$ws.Cells.Item(7, 10).Value = 297.5
$ws.Cells.Item(7, 12).Value = 297.5
$ws.Cells.Item(7, 14).Value = -523.5
$ws.Cells.Item(82, 8).Value = 26357.883
$ws.Cells.Item(82, 10).Value = 39997.5
$ws.Cells.Item(82, 12).Value = 39997.5
$ws.Cells.Item(82, 14).Value = -40763.5
$ws.Cells.Item(85, 8).Value = 26357.883
$ws.Cells.Item(85, 10).Value = 39997.5
$ws.Cells.Item(85, 12).Value = 39997.5
$ws.Cells.Item(85, 14).Value = -42649.5
$ws.Cells.Item(86, 8).Value = 2733.4075
$ws.Cells.Item(86, 9).Value = 2710.4
$ws.Cells.Item(86, 11).Value = 2710.4
$ws.Cells.Item(86, 13).Value = -1587.4
$ws.Cells.Item(89, 8).Value = 2733.4075
$ws.Cells.Item(89, 9).Value = 2710.4
$ws.Cells.Item(89, 11).Value = 13552
$ws.Cells.Item(89, 13).Value = -7936
$ws.Cells.Item(94, 8).Value = 1668.9
$ws.Cells.Item(94, 9).Value = 1763.1177
$ws.Cells.Item(94, 11).Value = 1763.1177
$ws.Cells.Item(94, 13).Value = -1312.1177
$ws.Cells.Item(107, 8).Value = 2325.5
$ws.Cells.Item(107, 9).Value = 988.25
$ws.Cells.Item(107, 11).Value = 988.25
$ws.Cells.Item(107, 13).Value = 931.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 13).Value = ""
$ws.Cells.Item(16, 8).Value = 8101.7144
$ws.Cells.Item(16, 10).Value = 9802.4
$ws.Cells.Item(16, 12).Value = 9802.4
$ws.Cells.Item(16, 14).Value = -10376.4
$ws.Cells.Item(31, 8).Value = 3905
$ws.Cells.Item(31, 9).Value = 2823.75
$ws.Cells.Item(31, 10).Value = 6067.5
$ws.Cells.Item(31, 11).Value = 2823.75
$ws.Cells.Item(31, 12).Value = 6067.5
$ws.Cells.Item(31, 13).Value = -2528.75
$ws.Cells.Item(31, 14).Value = -6657.5
$ws.Cells.Item(33, 8).Value = 280.14285
$ws.Cells.Item(33, 9).Value = 280.14285
$ws.Cells.Item(33, 11).Value = 280.14285
$ws.Cells.Item(33, 13).Value = 98.85714999999999
$ws.Cells.Item(34, 8).Value = 3905
$ws.Cells.Item(34, 9).Value = 2823.75
$ws.Cells.Item(34, 10).Value = 6067.5
$ws.Cells.Item(34, 11).Value = 2823.75
$ws.Cells.Item(34, 12).Value = 6067.5
$ws.Cells.Item(34, 13).Value = -2621.75
$ws.Cells.Item(34, 14).Value = -6471.5
$ws.Cells.Item(42, 8).Value = 9000
$ws.Cells.Item(42, 9).Value = 9000
$ws.Cells.Item(42, 11).Value = 9000
$ws.Cells.Item(42, 13).Value = -8407
$ws.Cells.Item(44, 8).Value = 7733
$ws.Cells.Item(44, 9).Value = 7733
$ws.Cells.Item(44, 11).Value = 7733
$ws.Cells.Item(44, 13).Value = -7291
$ws.Cells.Item(55, 8).Value = 12000
$ws.Cells.Item(55, 9).Value = 12000
$ws.Cells.Item(55, 11).Value = 12000
$ws.Cells.Item(55, 13).Value = -11685
$ws.Cells.Item(113, 8).Value = 8101.7144
$ws.Cells.Item(113, 10).Value = 9802.4
$ws.Cells.Item(113, 12).Value = 9802.4
$ws.Cells.Item(113, 14).Value = -14142.4
$ws.Cells.Item(122, 8).Value = 3112.4119
$ws.Cells.Item(122, 9).Value = 3742.75
$ws.Cells.Item(122, 11).Value = 11228.25
$ws.Cells.Item(122, 13).Value = -8778.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(29, 8).Value = 75.8
$ws.Cells.Item(29, 9).Value = 1
$ws.Cells.Item(29, 10).Value = 79.73684
$ws.Cells.Item(29, 11).Value = 3
$ws.Cells.Item(29, 12).Value = 239.21052
$ws.Cells.Item(29, 13).Value = 274
$ws.Cells.Item(29, 14).Value = -793.21052
$ws.Cells.Item(56, 8).Value = 18800.834
$ws.Cells.Item(56, 9).Value = 18800.834
$ws.Cells.Item(56, 11).Value = 18800.834
$ws.Cells.Item(56, 13).Value = -18270.834
$ws.Cells.Item(98, 8).Value = 388.5
$ws.Cells.Item(98, 9).Value = 388.5
$ws.Cells.Item(98, 11).Value = 1165.5
$ws.Cells.Item(98, 13).Value = 332.5
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 13).Value = ""
$ws.Cells.Item(134, 8).Value = 999.5
$ws.Cells.Item(134, 9).Value = 999.5
$ws.Cells.Item(134, 11).Value = 2998.5
$ws.Cells.Item(134, 13).Value = 2071.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8000
$ws.Cells.Item(70, 9).Value = 8000
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 8000
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = -7730
$ws.Cells.Item(70, 14).Value = ""
$ws.Cells.Item(73, 8).Value = 8000
$ws.Cells.Item(73, 9).Value = 8000
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 8000
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = -7064
$ws.Cells.Item(73, 14).Value = ""
$ws.Cells.Item(97, 8).Value = 731.94446
$ws.Cells.Item(97, 9).Value = 757.35297
$ws.Cells.Item(97, 10).Value = 300
$ws.Cells.Item(97, 11).Value = 757.35297
$ws.Cells.Item(97, 12).Value = 300
$ws.Cells.Item(97, 13).Value = -261.35297
$ws.Cells.Item(97, 14).Value = -1292
$ws.Cells.Item(126, 8).Value = 2256.125
$ws.Cells.Item(126, 9).Value = 2256.125
$ws.Cells.Item(126, 11).Value = 6768.375
$ws.Cells.Item(126, 13).Value = -4298.375
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).Value = ""
$ws.Cells.Item(136, 8).Value = 47662
$ws.Cells.Item(136, 10).Value = 47662
$ws.Cells.Item(136, 12).Value = 142986
$ws.Cells.Item(136, 14).Value = -148086

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5200.6
$ws.Cells.Item(7, 9).Value = 5000.75
$ws.Cells.Item(7, 11).Value = 5000.75
$ws.Cells.Item(7, 13).Value = -4888.75
$ws.Cells.Item(126, 8).Value = 5200.6
$ws.Cells.Item(126, 9).Value = 5000.75
$ws.Cells.Item(126, 11).Value = 15002.25
$ws.Cells.Item(126, 13).Value = -12532.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 499.25
$ws.Cells.Item(81, 9).Value = 499.25
$ws.Cells.Item(81, 11).Value = 998.5
$ws.Cells.Item(81, 13).Value = 62.5
$ws.Cells.Item(84, 8).Value = 499.25
$ws.Cells.Item(84, 9).Value = 499.25
$ws.Cells.Item(84, 11).Value = 4992.5
$ws.Cells.Item(84, 13).Value = 311.5
$ws.Cells.Item(100, 8).Value = 731.5
$ws.Cells.Item(100, 9).Value = 692.25
$ws.Cells.Item(100, 11).Value = 1384.5
$ws.Cells.Item(100, 13).Value = -843.5
$ws.Cells.Item(122, 8).Value = 1258.4
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 14).Value = ""
$ws.Cells.Item(123, 8).Value = 69997
$ws.Cells.Item(123, 10).Value = 69997
$ws.Cells.Item(123, 12).Value = 69997
$ws.Cells.Item(123, 14).Value = -79797
$ws.Cells.Item(132, 8).Value = 1308.1666
$ws.Cells.Item(132, 9).Value = 1177.5555
$ws.Cells.Item(132, 11).Value = 3532.6665
$ws.Cells.Item(132, 13).Value = -1002.6665
$ws.Cells.Item(136, 8).Value = 2535.2942
$ws.Cells.Item(136, 9).Value = 2599.75
$ws.Cells.Item(136, 11).Value = 7799.25
$ws.Cells.Item(136, 13).Value = -5249.25
